$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title.
#    Target shape:
#      <w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>
#            <w:r><w:t>: Experience ... Play for free now.</w:t></w:r></w:p>
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleEnd = $titlePara.Range.Duplicate
$titleEnd.Collapse(0)
$titleEnd.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = $d.Styles.Item("Normal")

# Borrow the run layout (leading empty run + text run) from an existing
# body paragraph by copy/pasting it, then overwrite its text in place --
# this preserves the leading empty <w:r/> that a fresh InsertParagraphAfter
# paragraph would not have.
$template = $d.Paragraphs.Item(5)
$template.Range.Copy()
$metaPara.Range.Paste()

$boldLen = ("Meta description").Length
$boldRange = $d.Range($metaPara.Range.Start, $metaPara.Range.Start + $boldLen)
$boldSourceText = $boldRange.Text
$boldRange.Bold = 1
$boldRange.Find.Execute($boldSourceText, $false, $false, $false, $false, $false, $true, 1, $false, "Meta description", 2)

$restStart = $metaPara.Range.Start + $boldLen
$restEnd = $metaPara.Range.End - 1
$restRange = $d.Range($restStart, $restEnd)
$restSourceText = $restRange.Text
$restRange.Find.Execute($restSourceText, $false, $false, $false, $false, $false, $true, 1, $false, ": Experience the underwater world of Atlantis Queen with exciting bonus features and win big. Play for free now.", 2)

# ------------------------------------------------------------------
# 2) Near the end: drop the duplicated bold "Play Atlantis Queen Free ..."
#    paragraph, and replace the italic paragraph's text with the new
#    image-prompt copy (keeping its italic run / leading empty run).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupTitlePara = $null
for ($i = $count; $i -ge 1; $i--) {
    if ($dupTitlePara -ne $null) {
        continue
    }
    $cand = $d.Paragraphs.Item($i)
    $candText = $cand.Range.Text.TrimEnd([char]13)
    if ($i -ne 1 -and $candText -eq "Play Atlantis Queen Free - Exciting Features and Big Prizes") {
        $dupTitlePara = $cand
    }
}
$dupTitlePara.Range.Delete()

$count2 = $d.Paragraphs.Count
$italicPara = $d.Paragraphs.Item($count2)
$italicStart = $italicPara.Range.Start
$italicEnd = $italicPara.Range.End - 1
$italicTextRange = $d.Range($italicStart, $italicEnd)
$italicTextRange.Text = "Please create an engaging feature image that fits the theme of Atlantis Queen slot game. The image should be in a cartoon style and feature a happy Maya warrior wearing glasses. The Maya warrior should be diving into the ocean with a treasure chest in hand, surrounded by colorful sea creatures and Atlantis ruins in the background. The image should capture the excitement and adventure of the Atlantis Queen game, enticing players to dive in and discover the treasures of the underwater world."
